$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.884.72"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "2.679.57"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'597.98"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").Value = "'174.73"
$ws.Range("E6").Value = "  -4.09%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.522"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("D9").Value = "2.677.88"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  -5.90%  "
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").Value = "'0.356"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "'4.98"
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("D14").Value = "3.170.58"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "71.935.91"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000184"
$ws.Range("E16").Value = "  -5.79%  "
$ws.Range("D17").Value = "'26.16"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").Value = "2.690.04"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "'12.19"
$ws.Range("E19").Value = "  +4.45%  "
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("D21").Value = "'370.52"
$ws.Range("E21").Value = "  -4.18%  "
$ws.Range("D22").Value = "'4.16"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("D24").Value = "'72.08"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'4.32"
$ws.Range("E26").Value = "  -3.60%  "
$ws.Range("D27").Value = "'9.75"
$ws.Range("E27").Value = "  -2.95%  "
$ws.Range("D28").Value = "2.818.53"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "0.0₃0969"
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("D31").Value = "'8.03"
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("D32").Value = "'501.92"
$ws.Range("E32").Value = "  -8.87%  "
$ws.Range("D33").Value = "'1.29"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("D34").Value = "'1.82"
$ws.Range("E34").Value = "  -2.19%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "'163.51"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").Value = "'19.54"
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("D38").Value = "'19.09"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").Value = "'1.37"
$ws.Range("E39").Value = "  -4.08%  "
$ws.Range("E40").Value = "  -5.24%  "
$ws.Range("E41").Value = "  -5.26%  "
$ws.Range("D43").Value = "'5.00"
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("D44").Value = "'2.55"
$ws.Range("E44").Value = "  -3.69%  "
$ws.Range("D45").Value = "'0.332"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'39.52"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'156.22"
$ws.Range("E47").Value = "  +2.03%  "
$ws.Range("D48").Value = "'0.563"
$ws.Range("E48").Value = "  +3.23%  "
$ws.Range("D49").Value = "'3.71"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").Value = "'1.74"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").Value = "'0.0755"
$ws.Range("E51").Value = "  -1.24%  "
